# Update ppt diagram files
#
# 1) Refresh the cached "datetimeFigureOut" date placeholder text
#    (12/15/2018 -> 4/1/2019) on the slide master and every slide layout.
# 2) Rename two shapes on slide 1: PersonListPanel -> EntryListPanel,
#    PersonCard -> EntryCard.

$p = $ppt.ActivePresentation

$oldDate = "12/15/2018"
$newDate = "4/1/2019"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($container) {
    for ($shapeIdx = 1; $shapeIdx -le $container.Shapes.Count; $shapeIdx++) {
        $sh = $container.Shapes.Item($shapeIdx)
        if ($sh.Type -eq 14) {
            if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $tr = $sh.TextFrame.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DatePlaceholder $master

# Every slide layout off the master
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout
}

# Slide 1 shape renames
$slide = $p.Slides.Item(1)
for ($si = 1; $si -le $slide.Shapes.Count; $si++) {
    $sh = $slide.Shapes.Item($si)
    if ($sh.HasTextFrame) {
        $txt = $sh.TextFrame.TextRange.Text
        if ($txt -eq "PersonListPanel") {
            $sh.TextFrame.TextRange.Text = "EntryListPanel"
        } elseif ($txt -eq "PersonCard") {
            $sh.TextFrame.TextRange.Text = "EntryCard"
        }
    }
}
